$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move 1.0 "At Work" hours to "Sick Leave" for Jan 2 and Jan 3 (rows 12 and 13)
$ws.Range("C12").Value = ""
$ws.Range("E12").Value = "1.0"

$ws.Range("C13").Value = ""
$ws.Range("E13").Value = "1.0"

# Update totals row (row 44): At Work 20.0 -> 18.0, Sick Leave "-" -> 2.0
$ws.Range("C44").Value = "18.0"
$ws.Range("E44").Value = "2.0"

# Update date label in row 48 (B48)
$ws.Range("B48").Value = "09 - February - 2025"
